# "worked on line balance"
# Adds a new "end time" (שעת_סיום) column (G) to Sheet1, mirroring the
# existing "start time" (שעת_התחלה) column F: same header styling
# (bold Arial, thin left/right border, centered), and the per-employee
# end-time values for the first 8 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell G1 ------------------------------------------------
$g1 = $ws.Range("G1")
$g1.Value = "שעת_סיום"

# Match F1's look: bold Arial 11, thin border on left+right only,
# centered both horizontally and vertically.
$g1.Font.Name = "Arial"
$g1.Font.Size = 11
$g1.Font.Bold = $true
$g1.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$g1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$g1.HorizontalAlignment = -4108     # xlCenter
$g1.VerticalAlignment = -4108       # xlCenter

# --- New data values G2:G9 ---------------------------------------------
$ws.Range("G2").Value = 23
$ws.Range("G3").Value = 18
$ws.Range("G4").Value = 23
$ws.Range("G5").Value = 23
$ws.Range("G6").Value = 23
$ws.Range("G7").Value = 18
$ws.Range("G8").Value = 23
$ws.Range("G9").Value = 18

# --- View tweaks (zoom + active selection moved to the new column) -----
$excel.ActiveWindow.Zoom = 160
$ws.Range("G9").Select() | Out-Null
